# Apply weekly Fruit/Vegetable data refresh: rows 2-8 in columns D, M, N, O, P, Q, R, S, T
# are reshuffled (the underlying daily records got re-sorted/re-dated), per the diff.
# This corresponds to a cyclic permutation where the "new" row N takes the values that
# used to live in row mapping[N] (see analysis of the unified diff).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that move together as a record's payload for rows 2-8.
$cols = @("D", "M", "N", "O", "P", "Q", "R", "S", "T")

# Mapping: new row number -> old row number whose values it should receive.
$rowMap = @{
    2 = 5
    3 = 6
    4 = 3
    5 = 7
    6 = 8
    7 = 4
    8 = 2
}

# Snapshot the original values for rows 2-8 before any writes, since several
# destinations and sources overlap (it's a single 7-cycle permutation).
$snapshot = @{}
foreach ($r in 2..8) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowVals
}

# Now write each destination row using the snapshot of its source row.
foreach ($newRow in $rowMap.Keys) {
    $oldRow = $rowMap[$newRow]
    $src = $snapshot[$oldRow]
    foreach ($c in $cols) {
        $ws.Range("$c$newRow").Value2 = $src[$c]
    }
}
